{"js": "// Update the date line and all 25 two-digit multiplication equations in\n// the table to their new values (matches the canonical OOXML diff).\nconst replacements = [\n  [\"2026-01-09 Friday\", \"2026-01-10 Saturday\"],\n  [\"68\u00d764=4352\", \"19\u00d777=1463\"],\n  [\"47\u00d713=611\", \"67\u00d772=4824\"],\n  [\"41\u00d733=1353\", \"94\u00d724=2256\"],\n  [\"53\u00d756=2968\", \"95\u00d784=7980\"],\n  [\"43\u00d799=4257\", \"83\u00d794=7802\"],\n  [\"76\u00d799=7524\", \"13\u00d783=1079\"],\n  [\"58\u00d782=4756\", \"92\u00d795=8740\"],\n  [\"47\u00d756=2632\", \"12\u00d792=1104\"],\n  [\"93\u00d751=4743\", \"22\u00d787=1914\"],\n  [\"31\u00d761=1891\", \"25\u00d759=1475\"],\n  [\"37\u00d757=2109\", \"50\u00d767=3350\"],\n  [\"80\u00d748=3840\", \"44\u00d799=4356\"],\n  [\"32\u00d714=448\", \"64\u00d726=1664\"],\n  [\"38\u00d789=3382\", \"81\u00d750=4050\"],\n  [\"74\u00d732=2368\", \"27\u00d780=2160\"],\n  [\"95\u00d739=3705\", \"13\u00d778=1014\"],\n  [\"36\u00d778=2808\", \"46\u00d759=2714\"],\n  [\"66\u00d744=2904\", \"55\u00d720=1100\"],\n  [\"36\u00d759=2124\", \"52\u00d730=1560\"],\n  [\"17\u00d776=1292\", \"13\u00d765=845\"],\n  [\"83\u00d720=1660\", \"70\u00d745=3150\"],\n  [\"48\u00d760=2880\", \"56\u00d749=2744\"],\n  [\"54\u00d737=1998\", \"20\u00d768=1360\"],\n  [\"24\u00d753=1272\", \"74\u00d715=1110\"],\n  [\"54\u00d773=3942\", \"34\u00d766=2244\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and all 25 two-digit multiplication equations in\n# the table to their new values (matches the canonical OOXML diff).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2026-01-09 Friday\", \"2026-01-10 Saturday\"),\n    @(\"68\u00d764=4352\", \"19\u00d777=1463\"),\n    @(\"47\u00d713=611\", \"67\u00d772=4824\"),\n    @(\"41\u00d733=1353\", \"94\u00d724=2256\"),\n    @(\"53\u00d756=2968\", \"95\u00d784=7980\"),\n    @(\"43\u00d799=4257\", \"83\u00d794=7802\"),\n    @(\"76\u00d799=7524\", \"13\u00d783=1079\"),\n    @(\"58\u00d782=4756\", \"92\u00d795=8740\"),\n    @(\"47\u00d756=2632\", \"12\u00d792=1104\"),\n    @(\"93\u00d751=4743\", \"22\u00d787=1914\"),\n    @(\"31\u00d761=1891\", \"25\u00d759=1475\"),\n    @(\"37\u00d757=2109\", \"50\u00d767=3350\"),\n    @(\"80\u00d748=3840\", \"44\u00d799=4356\"),\n    @(\"32\u00d714=448\", \"64\u00d726=1664\"),\n    @(\"38\u00d789=3382\", \"81\u00d750=4050\"),\n    @(\"74\u00d732=2368\", \"27\u00d780=2160\"),\n    @(\"95\u00d739=3705\", \"13\u00d778=1014\"),\n    @(\"36\u00d778=2808\", \"46\u00d759=2714\"),\n    @(\"66\u00d744=2904\", \"55\u00d720=1100\"),\n    @(\"36\u00d759=2124\", \"52\u00d730=1560\"),\n    @(\"17\u00d776=1292\", \"13\u00d765=845\"),\n    @(\"83\u00d720=1660\", \"70\u00d745=3150\"),\n    @(\"48\u00d760=2880\", \"56\u00d749=2744\"),\n    @(\"54\u00d737=1998\", \"20\u00d768=1360\"),\n    @(\"24\u00d753=1272\", \"74\u00d715=1110\"),\n    @(\"54\u00d773=3942\", \"34\u00d766=2244\")\n)\n\nforeach ($pair in $replacements) {\n    $find = $pair[0]\n    $replace = $pair[1]\n\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Text = $find\n    $range.Find.Replacement.Text = $replace\n    $range.Find.Forward = $true\n    $range.Find.Wrap = 1\n    $range.Find.MatchCase = $true\n    $range.Find.MatchWholeWord = $false\n    $range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)\n}\n"}
